$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of this data block (rows 205-207),
# pushing the existing rows 205-213 down to 208-216.
$ws.Rows("205:207").Insert()

# ---- Row 205 (new) ----
$ws.Range("A205").Value = 10
$ws.Range("B205").Value = "Vega Modelo de Temuco"
$ws.Range("C205").Value = "La Araucanía"
$ws.Range("D205").Value = 44578
$ws.Range("E205").Value = 9
$ws.Range("F205").Value = "Fruta"
$ws.Range("G205").Value = 100103
$ws.Range("H205").Value = "Frutos de hueso (carozo)"
$ws.Range("I205").Value = 100103004
$ws.Range("J205").Value = "Durazno"
$ws.Range("K205").Value = "Kurakata"
$ws.Range("L205").Value = "Especial"
$ws.Range("M205").Value = 95
$ws.Range("N205").Value = 20000
$ws.Range("O205").Value = 20000
$ws.Range("P205").Value = 20000
$ws.Range("Q205").Value = "$/bandeja 18 kilos granel"
$ws.Range("R205").Value = "Región de O'Higgins"
$ws.Range("S205").Value = 1111
$ws.Range("T205").Value = 18

# ---- Row 206 (new) ----
$ws.Range("A206").Value = 10
$ws.Range("B206").Value = "Vega Modelo de Temuco"
$ws.Range("C206").Value = "La Araucanía"
$ws.Range("D206").Value = 44578
$ws.Range("E206").Value = 9
$ws.Range("F206").Value = "Fruta"
$ws.Range("G206").Value = 100103
$ws.Range("H206").Value = "Frutos de hueso (carozo)"
$ws.Range("I206").Value = 100103004
$ws.Range("J206").Value = "Durazno"
$ws.Range("K206").Value = "Kurakata"
$ws.Range("L206").Value = "Primera"
$ws.Range("M206").Value = 110
$ws.Range("N206").Value = 18000
$ws.Range("O206").Value = 18000
$ws.Range("P206").Value = 18000
$ws.Range("Q206").Value = "$/bandeja 18 kilos granel"
$ws.Range("R206").Value = "Región de O'Higgins"
$ws.Range("S206").Value = 1000
$ws.Range("T206").Value = 18

# ---- Row 207 (new) ----
$ws.Range("A207").Value = 10
$ws.Range("B207").Value = "Vega Modelo de Temuco"
$ws.Range("C207").Value = "La Araucanía"
$ws.Range("D207").Value = 44578
$ws.Range("E207").Value = 9
$ws.Range("F207").Value = "Fruta"
$ws.Range("G207").Value = 100103
$ws.Range("H207").Value = "Frutos de hueso (carozo)"
$ws.Range("I207").Value = 100103004
$ws.Range("J207").Value = "Durazno"
$ws.Range("K207").Value = "Kurakata"
$ws.Range("L207").Value = "Segunda"
$ws.Range("M207").Value = 95
$ws.Range("N207").Value = 14000
$ws.Range("O207").Value = 14000
$ws.Range("P207").Value = 14000
$ws.Range("Q207").Value = "$/bandeja 18 kilos granel"
$ws.Range("R207").Value = "Región de O'Higgins"
$ws.Range("S207").Value = 778
$ws.Range("T207").Value = 18

Write-Host "Done. UsedRange rows:" $ws.UsedRange.Rows.Count
